$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Column I (Pre_ISI) values per row
$ws.Range("I2").Value = 1.75
$ws.Range("I3").Value = 2.9399999999999995
$ws.Range("I4").Value = 2.2200000000000006
$ws.Range("I5").Value = 1.67
$ws.Range("I6").Value = 2.0099999999999998
$ws.Range("I7").Value = 2.79
$ws.Range("I8").Value = 2.2300000000000004
$ws.Range("I9").Value = 2.76
$ws.Range("I10").Value = 1.58
$ws.Range("I11").Value = 2.9700000000000006
$ws.Range("I12").Value = 2.3100000000000005
$ws.Range("I13").Value = 2.9800000000000004
$ws.Range("I14").Value = 2.95
$ws.Range("I15").Value = 2.76
$ws.Range("I16").Value = 2.5700000000000003
$ws.Range("I17").Value = 1.6600000000000001
$ws.Range("I18").Value = 2.99
$ws.Range("I19").Value = 2.95
$ws.Range("I20").Value = 2.79
$ws.Range("I21").Value = 2.96
$ws.Range("I22").Value = 2.0499999999999998
$ws.Range("I23").Value = 2.3100000000000005
$ws.Range("I24").Value = 2.08
$ws.Range("I25").Value = 2.6500000000000004
$ws.Range("I26").Value = 1.52
$ws.Range("I27").Value = 2.42
$ws.Range("I28").Value = 2.1099999999999994
$ws.Range("I29").Value = 1.6400000000000001
$ws.Range("I30").Value = 2.3599999999999994
$ws.Range("I31").Value = 2.8
$ws.Range("I32").Value = 2.3599999999999994
$ws.Range("I33").Value = 2.4000000000000004
$ws.Range("I34").Value = 1.98
$ws.Range("I35").Value = 1.5699999999999998
$ws.Range("I36").Value = 2.99
$ws.Range("I37").Value = 2.0700000000000003
$ws.Range("I38").Value = 2.7200000000000006
$ws.Range("I39").Value = 2.6899999999999995
$ws.Range("I40").Value = 2.3900000000000006
$ws.Range("I41").Value = 2.3100000000000005
$ws.Range("I42").Value = 2.9000000000000004
$ws.Range("I43").Value = 1.77
$ws.Range("I44").Value = 2.5099999999999998
$ws.Range("I45").Value = 1.6800000000000002
$ws.Range("I46").Value = 1.58
$ws.Range("I47").Value = 2.25
$ws.Range("I48").Value = 2.58
$ws.Range("I49").Value = 1.6400000000000001
$ws.Range("I50").Value = 2.2400000000000002
$ws.Range("I51").Value = 2.4000000000000004
$ws.Range("I52").Value = 2.63
$ws.Range("I53").Value = 1.8599999999999999
$ws.Range("I54").Value = 2.2300000000000004
$ws.Range("I55").Value = 1.81
$ws.Range("I56").Value = 1.8599999999999999
$ws.Range("I57").Value = 2.92
$ws.Range("I58").Value = 2.83
$ws.Range("I59").Value = 2.6099999999999994
$ws.Range("I60").Value = 2.1799999999999997
$ws.Range("I61").Value = 2.34
$ws.Range("I62").Value = 2.92
$ws.Range("I63").Value = 1.94
$ws.Range("I64").Value = 1.9700000000000002
$ws.Range("I65").Value = 2.9700000000000006
$ws.Range("I66").Value = 1.8599999999999999
$ws.Range("I67").Value = 2.3200000000000003
$ws.Range("I68").Value = 2.74
$ws.Range("I69").Value = 1.5699999999999998
$ws.Range("I70").Value = 2.2999999999999998
$ws.Range("I71").Value = 2.0499999999999998
$ws.Range("I72").Value = 2.9800000000000004
$ws.Range("I73").Value = 1.9500000000000002
$ws.Range("I74").Value = 2.79
$ws.Range("I75").Value = 2.13
$ws.Range("I76").Value = 2.12
$ws.Range("I77").Value = 2.7
$ws.Range("I78").Value = 2.7
$ws.Range("I79").Value = 2.13
$ws.Range("I80").Value = 2.2400000000000002
$ws.Range("I81").Value = 1.77
$ws.Range("I82").Value = 2.5999999999999996
$ws.Range("I83").Value = 1.5
$ws.Range("I84").Value = 2.8
$ws.Range("I85").Value = 1.8399999999999999
$ws.Range("I86").Value = 1.5899999999999999
$ws.Range("I87").Value = 2.9000000000000004
$ws.Range("I88").Value = 2.2000000000000002
$ws.Range("I89").Value = 1.9300000000000002
$ws.Range("I90").Value = 2.99
$ws.Range("I91").Value = 2.8200000000000003
$ws.Range("I92").Value = 2.88
$ws.Range("I93").Value = 1.73
$ws.Range("I94").Value = 2.92
$ws.Range("I95").Value = 1.63
$ws.Range("I96").Value = 1.83
$ws.Range("I97").Value = 1.83
$ws.Range("I98").Value = 2.42
$ws.Range("I99").Value = 1.77
$ws.Range("I100").Value = 2.21
$ws.Range("I101").Value = 2.1799999999999997
$ws.Range("I102").Value = 2.3100000000000005
$ws.Range("I103").Value = 1.6099999999999999
$ws.Range("I104").Value = 2.5299999999999994
$ws.Range("I105").Value = 2.2000000000000002
$ws.Range("I106").Value = 1.54
$ws.Range("I107").Value = 2.08
$ws.Range("I108").Value = 2.0299999999999994
$ws.Range("I109").Value = 1.8199999999999998

# Update selected cell to match final workbook state
$ws.Range("L14").Select()
